$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update timestamp column Z (col 26) for rows 2-48 with new timestamp values
$timestamps = @{
    2 = "2025-10-17T07:09:26.433716"
    3 = "2025-10-17T07:09:26.433716"
    4 = "2025-10-17T07:09:26.434716"
    5 = "2025-10-17T07:09:26.434716"
    6 = "2025-10-17T07:09:26.434716"
    7 = "2025-10-17T07:09:26.434716"
    8 = "2025-10-17T07:09:26.434716"
    9 = "2025-10-17T07:09:26.434716"
    10 = "2025-10-17T07:09:26.434716"
    11 = "2025-10-17T07:09:26.434716"
    12 = "2025-10-17T07:09:26.434716"
    13 = "2025-10-17T07:09:26.434716"
    14 = "2025-10-17T07:09:26.435714"
    15 = "2025-10-17T07:09:26.435714"
    16 = "2025-10-17T07:09:26.491276"
    17 = "2025-10-17T07:09:26.491789"
    18 = "2025-10-17T07:09:26.491789"
    19 = "2025-10-17T07:09:26.491789"
    20 = "2025-10-17T07:09:26.491789"
    21 = "2025-10-17T07:09:26.491789"
    22 = "2025-10-17T07:09:26.492804"
    23 = "2025-10-17T07:09:26.492804"
    24 = "2025-10-17T07:09:26.492804"
    25 = "2025-10-17T07:09:26.492804"
    26 = "2025-10-17T07:09:26.557762"
    27 = "2025-10-17T07:09:26.557762"
    28 = "2025-10-17T07:09:26.557762"
    29 = "2025-10-17T07:09:26.557762"
    30 = "2025-10-17T07:09:26.558761"
    31 = "2025-10-17T07:09:26.558761"
    32 = "2025-10-17T07:09:26.558761"
    33 = "2025-10-17T07:09:26.559761"
    34 = "2025-10-17T07:09:26.559761"
    35 = "2025-10-17T07:09:26.559761"
    36 = "2025-10-17T07:09:26.560764"
    37 = "2025-10-17T07:09:26.560764"
    38 = "2025-10-17T07:09:26.560764"
    39 = "2025-10-17T07:09:26.560764"
    40 = "2025-10-17T07:09:26.561764"
    41 = "2025-10-17T07:09:26.561764"
    42 = "2025-10-17T07:09:26.561764"
    43 = "2025-10-17T07:09:26.561764"
    44 = "2025-10-17T07:09:26.561764"
    45 = "2025-10-17T07:09:26.562763"
    46 = "2025-10-17T07:09:26.562763"
    47 = "2025-10-17T07:09:26.562763"
    48 = "2025-10-17T07:09:26.562763"
}

foreach ($row in $timestamps.Keys) {
    $ws.Cells.Item($row, 26).Value = $timestamps[$row]
}
